$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 3 entirely; rows below (4-7) shift up to (3-6).
$ws.Rows.Item(3).Delete()
